# Updates the premier-league 2023-2024 match list:
#  1) Re-shuffles the F:V (match detail) columns across several existing
#     rows back into their correct chronological slot (A:E - index/country/
#     tournament/season/date - stay put; only the match info moves).
#  2) Appends 4 new match rows (113-116) at the end of the table, copying
#     the row-112 formatting down so the new rows match the existing style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Row permutation -----------------------------------------------
# new row number -> old row number whose F:V contents it should receive.
$mapping = @{
    42 = 44; 43 = 46; 44 = 45; 45 = 43; 46 = 42;
    51 = 52; 52 = 51;
    57 = 58; 58 = 57;
    74 = 75; 75 = 76; 76 = 74;
    83 = 87; 84 = 86; 85 = 83; 86 = 85; 87 = 84;
    98 = 99; 99 = 100; 100 = 98;
    105 = 106; 106 = 105;
}

# Snapshot the F:V values of every row that participates in a permutation
# BEFORE any writes happen (writes happen in place, so we must not read
# from a row we may have already overwritten).
$snapshot = @{}
foreach ($row in $mapping.Values) {
    if (-not $snapshot.ContainsKey($row)) {
        $snapshot[$row] = $ws.Range("F$row`:V$row").Value2
    }
}

foreach ($row in $mapping.Keys) {
    $srcRow = $mapping[$row]
    $ws.Range("F$row`:V$row").Value2 = $snapshot[$srcRow]
}

# --- 2) Append new rows 113-116 ----------------------------------------
# Copy the formatting of the last existing row (112) down onto the new
# rows so they inherit the same styles (bold/boxed index column, date
# number format, etc.), then overwrite with the real values.
$ws.Range("A112:V112").Copy($ws.Range("A113:V116"))

$newRows = @(
    @{ Row = 113; A = 112; F = "Crystal Palace"; G = 2; H = "Everton"; I = 3;
       E = 45241.66666666666;
       J = 2.09;  K = "28/10/2023 20:02"; L = 2.71;  M = "11/11/2023 15:58";
       N = 3.34;  O = "28/10/2023 20:02"; P = 3.12;  Q = "11/11/2023 15:56";
       R = 3.8;   S = "28/10/2023 20:02"; T = 2.96;  U = "11/11/2023 15:59";
       V = "https://www.betexplorer.com/football/england/premier-league/crystal-palace-everton/rZrW8iVi/" },

    @{ Row = 114; A = 113; F = "Arsenal"; G = 3; H = "Burnley"; I = 1;
       E = 45241.66666666666;
       J = 1.24;  K = "29/10/2023 11:22"; L = 1.19;  M = "11/11/2023 15:13";
       N = 6.5;   O = "29/10/2023 11:22"; P = 7.49;  Q = "11/11/2023 15:36";
       R = 12.06; S = "29/10/2023 11:22"; T = 16.59; U = "11/11/2023 15:36";
       V = "https://www.betexplorer.com/football/england/premier-league/arsenal-burnley/ncYLjAFN/" },

    @{ Row = 115; A = 114; F = "Manchester Utd"; G = 1; H = "Luton"; I = 0;
       E = 45241.66666666666;
       J = 1.2;   K = "29/10/2023 11:21"; L = 1.29;  M = "11/11/2023 15:58";
       N = 7.39;  O = "29/10/2023 11:21"; P = 6.14;  Q = "11/11/2023 15:59";
       R = 13.36; S = "29/10/2023 11:21"; T = 10.06; U = "11/11/2023 15:59";
       V = "https://www.betexplorer.com/football/england/premier-league/manchester-united-luton/4lXv7Va4/" },

    @{ Row = 116; A = 115; F = "Bournemouth"; G = 2; H = "Newcastle"; I = 0;
       E = 45241.77083333334;
       J = 4.92;  K = "29/10/2023 00:02"; L = 3.9;   M = "11/11/2023 18:28";
       N = 4.18;  O = "29/10/2023 00:02"; P = 3.76;  Q = "11/11/2023 18:28";
       R = 1.67;  S = "29/10/2023 00:02"; T = 1.97;  U = "11/11/2023 18:28";
       V = "https://www.betexplorer.com/football/england/premier-league/bournemouth-newcastle-utd/86ZHij0H/" }
)

foreach ($nr in $newRows) {
    $row = $nr.Row
    $ws.Cells.Item($row, 1).Value2 = $nr.A          # A - Indice
    $ws.Cells.Item($row, 2).Value2 = "england"       # B - pais
    $ws.Cells.Item($row, 3).Value2 = "premier-league" # C - torneio
    $ws.Cells.Item($row, 4).Value2 = "2023-2024"     # D - temporada
    $ws.Cells.Item($row, 5).Value2 = $nr.E           # E - data_partida
    $ws.Cells.Item($row, 6).Value2 = $nr.F           # F - home
    $ws.Cells.Item($row, 7).Value2 = $nr.G           # G - home_ft_gols
    $ws.Cells.Item($row, 8).Value2 = $nr.H           # H - away
    $ws.Cells.Item($row, 9).Value2 = $nr.I           # I - away_ft_gols
    $ws.Cells.Item($row, 10).Value2 = $nr.J          # J - home_opening_odds
    $ws.Cells.Item($row, 11).Value2 = $nr.K          # K - home_opening_data_hora
    $ws.Cells.Item($row, 12).Value2 = $nr.L          # L - home_closing_odds
    $ws.Cells.Item($row, 13).Value2 = $nr.M          # M - home_closing_data_hora
    $ws.Cells.Item($row, 14).Value2 = $nr.N          # N - draw_opening_odds
    $ws.Cells.Item($row, 15).Value2 = $nr.O          # O - draw_opening_data_hora
    $ws.Cells.Item($row, 16).Value2 = $nr.P          # P - draw_closing_odds
    $ws.Cells.Item($row, 17).Value2 = $nr.Q          # Q - draw_closing_data_hora
    $ws.Cells.Item($row, 18).Value2 = $nr.R          # R - away_opening_odds
    $ws.Cells.Item($row, 19).Value2 = $nr.S          # S - away_opening_data_hora
    $ws.Cells.Item($row, 20).Value2 = $nr.T          # T - away_closing_odds
    $ws.Cells.Item($row, 21).Value2 = $nr.U          # U - away_closing_data_hora
    $ws.Cells.Item($row, 22).Value2 = $nr.V          # V - url_partida
}
